$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for columns C (date), D (Vinculos), E (Diferenca), F (Variacao %)
# representing the data shifted by one year block (8 rows per year) with a new 2023 block appended
$data = New-Object 'object[,]' 80,4
$data[0,0] = "31/12/2014"; $data[0,1] = 4894; $data[0,2] = 290; $data[0,3] = 6.3
$data[1,0] = "31/12/2014"; $data[1,1] = 48306; $data[1,2] = 1145; $data[1,3] = 2.43
$data[2,0] = "31/12/2014"; $data[2,1] = 6188; $data[2,2] = 203; $data[2,3] = 3.39
$data[3,0] = "31/12/2014"; $data[3,1] = 27674; $data[3,2] = -2198; $data[3,3] = -7.36
$data[4,0] = "31/12/2014"; $data[4,1] = 68655; $data[4,2] = 3161; $data[4,3] = 4.83
$data[5,0] = "31/12/2014"; $data[5,1] = 133255; $data[5,2] = 8999; $data[5,3] = 7.24
$data[6,0] = "31/12/2014"; $data[6,1] = 114287; $data[6,2] = -1695; $data[6,3] = -1.46
$data[7,0] = "31/12/2014"; $data[7,1] = 13764; $data[7,2] = 1343; $data[7,3] = 10.81
$data[8,0] = "31/12/2015"; $data[8,1] = 4296; $data[8,2] = -598; $data[8,3] = -12.22
$data[9,0] = "31/12/2015"; $data[9,1] = 48173; $data[9,2] = -133; $data[9,3] = -0.28
$data[10,0] = "31/12/2015"; $data[10,1] = 6076; $data[10,2] = -112; $data[10,3] = -1.81
$data[11,0] = "31/12/2015"; $data[11,1] = 24703; $data[11,2] = -2971; $data[11,3] = -10.74
$data[12,0] = "31/12/2015"; $data[12,1] = 67946; $data[12,2] = -709; $data[12,3] = -1.03
$data[13,0] = "31/12/2015"; $data[13,1] = 137164; $data[13,2] = 3909; $data[13,3] = 2.93
$data[14,0] = "31/12/2015"; $data[14,1] = 105628; $data[14,2] = -8659; $data[14,3] = -7.58
$data[15,0] = "31/12/2015"; $data[15,1] = 10982; $data[15,2] = -2782; $data[15,3] = -20.21
$data[16,0] = "31/12/2016"; $data[16,1] = 3595; $data[16,2] = -701; $data[16,3] = -16.32
$data[17,0] = "31/12/2016"; $data[17,1] = 43214; $data[17,2] = -4959; $data[17,3] = -10.29
$data[18,0] = "31/12/2016"; $data[18,1] = 4983; $data[18,2] = -1093; $data[18,3] = -17.99
$data[19,0] = "31/12/2016"; $data[19,1] = 17332; $data[19,2] = -7371; $data[19,3] = -29.84
$data[20,0] = "31/12/2016"; $data[20,1] = 65635; $data[20,2] = -2311; $data[20,3] = -3.4
$data[21,0] = "31/12/2016"; $data[21,1] = 133542; $data[21,2] = -3622; $data[21,3] = -2.64
$data[22,0] = "31/12/2016"; $data[22,1] = 103838; $data[22,2] = -1790; $data[22,3] = -1.69
$data[23,0] = "31/12/2016"; $data[23,1] = 10934; $data[23,2] = -48; $data[23,3] = -0.44
$data[24,0] = "31/12/2017"; $data[24,1] = 3048; $data[24,2] = -547; $data[24,3] = -15.22
$data[25,0] = "31/12/2017"; $data[25,1] = 42099; $data[25,2] = -1115; $data[25,3] = -2.58
$data[26,0] = "31/12/2017"; $data[26,1] = 5678; $data[26,2] = 695; $data[26,3] = 13.95
$data[27,0] = "31/12/2017"; $data[27,1] = 15447; $data[27,2] = -1885; $data[27,3] = -10.88
$data[28,0] = "31/12/2017"; $data[28,1] = 66061; $data[28,2] = 426; $data[28,3] = 0.65
$data[29,0] = "31/12/2017"; $data[29,1] = 132640; $data[29,2] = -902; $data[29,3] = -0.68
$data[30,0] = "31/12/2017"; $data[30,1] = 114675; $data[30,2] = 10837; $data[30,3] = 10.44
$data[31,0] = "31/12/2017"; $data[31,1] = 10814; $data[31,2] = -120; $data[31,3] = -1.1
$data[32,0] = "31/12/2018"; $data[32,1] = 2685; $data[32,2] = -363; $data[32,3] = -11.91
$data[33,0] = "31/12/2018"; $data[33,1] = 42246; $data[33,2] = 147; $data[33,3] = 0.35
$data[34,0] = "31/12/2018"; $data[34,1] = 5732; $data[34,2] = 54; $data[34,3] = 0.95
$data[35,0] = "31/12/2018"; $data[35,1] = 14643; $data[35,2] = -804; $data[35,3] = -5.2
$data[36,0] = "31/12/2018"; $data[36,1] = 64093; $data[36,2] = -1968; $data[36,3] = -2.98
$data[37,0] = "31/12/2018"; $data[37,1] = 137808; $data[37,2] = 5168; $data[37,3] = 3.9
$data[38,0] = "31/12/2018"; $data[38,1] = 111411; $data[38,2] = -3264; $data[38,3] = -2.85
$data[39,0] = "31/12/2018"; $data[39,1] = 10733; $data[39,2] = -81; $data[39,3] = -0.75
$data[40,0] = "31/12/2019"; $data[40,1] = 3162; $data[40,2] = 477; $data[40,3] = 17.77
$data[41,0] = "31/12/2019"; $data[41,1] = 38654; $data[41,2] = -3592; $data[41,3] = -8.5
$data[42,0] = "31/12/2019"; $data[42,1] = 6255; $data[42,2] = 523; $data[42,3] = 9.12
$data[43,0] = "31/12/2019"; $data[43,1] = 13339; $data[43,2] = -1304; $data[43,3] = -8.91
$data[44,0] = "31/12/2019"; $data[44,1] = 65263; $data[44,2] = 1170; $data[44,3] = 1.83
$data[45,0] = "31/12/2019"; $data[45,1] = 137457; $data[45,2] = -351; $data[45,3] = -0.25
$data[46,0] = "31/12/2019"; $data[46,1] = 78380; $data[46,2] = -33031; $data[46,3] = -29.65
$data[47,0] = "31/12/2019"; $data[47,1] = 10360; $data[47,2] = -373; $data[47,3] = -3.48
$data[48,0] = "31/12/2020"; $data[48,1] = 2405; $data[48,2] = -757; $data[48,3] = -23.94
$data[49,0] = "31/12/2020"; $data[49,1] = 41259; $data[49,2] = 2605; $data[49,3] = 6.74
$data[50,0] = "31/12/2020"; $data[50,1] = 6248; $data[50,2] = -7; $data[50,3] = -0.11
$data[51,0] = "31/12/2020"; $data[51,1] = 14998; $data[51,2] = 1659; $data[51,3] = 12.44
$data[52,0] = "31/12/2020"; $data[52,1] = 64082; $data[52,2] = -1181; $data[52,3] = -1.81
$data[53,0] = "31/12/2020"; $data[53,1] = 130608; $data[53,2] = -6849; $data[53,3] = -4.98
$data[54,0] = "31/12/2020"; $data[54,1] = 96772; $data[54,2] = 18392; $data[54,3] = 23.47
$data[55,0] = "31/12/2020"; $data[55,1] = 9682; $data[55,2] = -678; $data[55,3] = -6.54
$data[56,0] = "31/12/2021"; $data[56,1] = 2279; $data[56,2] = -126; $data[56,3] = -5.24
$data[57,0] = "31/12/2021"; $data[57,1] = 41915; $data[57,2] = 656; $data[57,3] = 1.59
$data[58,0] = "31/12/2021"; $data[58,1] = 6248; $data[58,2] = 0; $data[58,3] = 0
$data[59,0] = "31/12/2021"; $data[59,1] = 16996; $data[59,2] = 1998; $data[59,3] = 13.32
$data[60,0] = "31/12/2021"; $data[60,1] = 68351; $data[60,2] = 4269; $data[60,3] = 6.66
$data[61,0] = "31/12/2021"; $data[61,1] = 134574; $data[61,2] = 3966; $data[61,3] = 3.04
$data[62,0] = "31/12/2021"; $data[62,1] = 110115; $data[62,2] = 13343; $data[62,3] = 13.79
$data[63,0] = "31/12/2021"; $data[63,1] = 11592; $data[63,2] = 1910; $data[63,3] = 19.73
$data[64,0] = "31/12/2022"; $data[64,1] = 2866; $data[64,2] = 587; $data[64,3] = 25.76
$data[65,0] = "31/12/2022"; $data[65,1] = 41531; $data[65,2] = -384; $data[65,3] = -0.92
$data[66,0] = "31/12/2022"; $data[66,1] = 5898; $data[66,2] = -350; $data[66,3] = -5.6
$data[67,0] = "31/12/2022"; $data[67,1] = 21123; $data[67,2] = 4127; $data[67,3] = 24.28
$data[68,0] = "31/12/2022"; $data[68,1] = 73886; $data[68,2] = 5535; $data[68,3] = 8.1
$data[69,0] = "31/12/2022"; $data[69,1] = 147605; $data[69,2] = 13031; $data[69,3] = 9.68
$data[70,0] = "31/12/2022"; $data[70,1] = 110761; $data[70,2] = 646; $data[70,3] = 0.59
$data[71,0] = "31/12/2022"; $data[71,1] = 13103; $data[71,2] = 1511; $data[71,3] = 13.03
$data[72,0] = "31/12/2023"; $data[72,1] = 3143; $data[72,2] = 277; $data[72,3] = 9.67
$data[73,0] = "31/12/2023"; $data[73,1] = 45289; $data[73,2] = 3758; $data[73,3] = 9.05
$data[74,0] = "31/12/2023"; $data[74,1] = 7006; $data[74,2] = 1108; $data[74,3] = 18.79
$data[75,0] = "31/12/2023"; $data[75,1] = 23979; $data[75,2] = 2856; $data[75,3] = 13.52
$data[76,0] = "31/12/2023"; $data[76,1] = 76807; $data[76,2] = 2921; $data[76,3] = 3.95
$data[77,0] = "31/12/2023"; $data[77,1] = 160044; $data[77,2] = 12439; $data[77,3] = 8.43
$data[78,0] = "31/12/2023"; $data[78,1] = 122930; $data[78,2] = 12169; $data[78,3] = 10.99
$data[79,0] = "31/12/2023"; $data[79,1] = 12906; $data[79,2] = -197; $data[79,3] = -1.5

$ws.Range("C2:F81").Value = $data
